# ---------------------------------------------------------------------------
# Feature: allow different intra-problem sets per variable (GitHub issue #75)
# Update the "_set_products", "_set_technologies" and "_set_years" sheets,
# and add two brand-new sets: "_set_hours" and "_set_scenarios".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper: copy header style (bold/centered/bordered) from a known header
#     cell onto a destination cell, then write its text ------------------
function Set-Header($destRange, $text, $styleSourceRange) {
    $styleSourceRange.Copy()
    $destRange.PasteSpecial(-4122)   # xlPasteFormats
    $destRange.Value = $text
}

# ===========================================================================
# Sheet 1 : _set_products  -> drop the "goods" member (row 4)
# ===========================================================================
$wsProducts = $wb.Worksheets.Item(1)
$wsProducts.Range("A4").Value = $null
$wsProducts.Range("A3").Select() | Out-Null

# ===========================================================================
# Sheet 2 : _set_technologies -> replace technology members/categories
# ===========================================================================
$wsTech = $wb.Worksheets.Item(2)
$wsTech.Range("A2").Value = "housing"
$wsTech.Range("B2").Value = "detailed"
$wsTech.Range("A3").Value = "transport"
$wsTech.Range("B3").Value = "detailed"
$wsTech.Range("A4").Value = "service"
$wsTech.Range("B4").Value = "not-detailed"
$wsTech.Columns.Item(1).ColumnWidth = 11.79
$wsTech.Columns.Item(2).ColumnWidth = 14.36
$wsTech.Range("A3").Select() | Out-Null

# ===========================================================================
# Sheet 3 : _set_years -> shrink from 6 years (2025-2030) to 3 (2020-2022)
# ===========================================================================
$wsYears = $wb.Worksheets.Item(3)
$wsYears.Range("A2").Value = 2020
$wsYears.Range("B2").Value = "warmup"
$wsYears.Range("A3").Value = 2021
$wsYears.Range("B3").Value = "run"
$wsYears.Range("A4").Value = 2022
$wsYears.Range("B4").Value = "run"
$wsYears.Range("A5").Value = $null
$wsYears.Range("B5").Value = $null
$wsYears.Range("A6").Value = $null
$wsYears.Range("B6").Value = $null
$wsYears.Range("A7").Value = $null
$wsYears.Range("B7").Value = $null
$wsYears.Range("A3").Select() | Out-Null

# ===========================================================================
# Sheet 4 (new) : _set_hours
# ===========================================================================
$wsHours = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsYears)
$wsHours.Name = "_set_hours"
Set-Header $wsHours.Range("A1") "h_Names" $wsTech.Range("A1")
$wsHours.Range("A2").Value = "h1"
$wsHours.Range("A3").Value = "h2"
$wsHours.Range("A4").Value = "h3"
$wsHours.Range("A5").Value = "h4"
$wsHours.Range("A6").Value = "h5"
$wsHours.Range("A3").Select() | Out-Null

# ===========================================================================
# Sheet 5 (new) : _set_scenarios
# ===========================================================================
$wsScenarios = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsHours)
$wsScenarios.Name = "_set_scenarios"
Set-Header $wsScenarios.Range("A1") "s_Names" $wsTech.Range("A1")
$wsScenarios.Range("A2").Value = "step"
$wsScenarios.Range("A3").Value = "nze"
$wsScenarios.Range("A4").Value = "opt"
$wsScenarios.Range("G3").Select() | Out-Null
$wsScenarios.Activate() | Out-Null

Write-Host "Stage 2 done"
